$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-generated NET_STATION_ID values (column A) replace the previous
# "romN" placeholders with freshly generated station ids.
$ws.Range("A2").Value  = "Auto2020-02-06-6978"
$ws.Range("A3").Value  = "Auto2020-02-06-6709"
$ws.Range("A4").Value  = "Auto2020-02-06-8166"
$ws.Range("A5").Value  = "Auto2020-02-06-6877"
$ws.Range("A6").Value  = "Auto2020-02-06-4736"
$ws.Range("A7").Value  = "Auto2020-02-06-5273"
$ws.Range("A8").Value  = "Auto2020-02-06-3194"
$ws.Range("A9").Value  = "Auto2020-02-06-6478"
$ws.Range("A10").Value = "Auto2020-02-06-5578"
$ws.Range("A11").Value = "Auto2020-02-06-5701"
$ws.Range("A12").Value = "Auto2020-02-06-1499"
$ws.Range("A13").Value = "Auto2020-02-06-8967"
$ws.Range("A14").Value = "Auto2020-02-06-513"
$ws.Range("A15").Value = "Auto2020-02-06-796"

# ERRAND_DATE (column J) rewritten to the same date for each data row as
# part of the same pass.
$ws.Range("J2:J15").Value = "2019-12-06"

# Leave the active selection on A5 (previously K3:K15).
[void]$ws.Range("A5").Select()
